$wb = $excel.ActiveWorkbook

# Rename "Sheet2" -> "单元表2"
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "单元表2"

# Make the renamed sheet the active one (tabSelected/activeTab)
$ws2.Activate()
